$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"

$ws.Range("D2").Value = "12/29/2025"
$ws.Range("D3").Value = "12/29/2025"
